# Reorder PCB congeners in the data table.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row order (header stays in row 1): sumPCB moves to the top, then the
# "CB2x"/"CB3x"/"CB4x"/"CB5x" congeners, followed by the remaining CB1xx
# congeners in their original relative order.
$data = @(
    @("sumPCB", 0.2,  1.22, 6.49),
    @("CB28",   0.06, 0.6,  2.64),
    @("CB31",   0.01, 0.31, 1.76),
    @("CB44",   0.09, 0.92, 6.22),
    @("CB49",   0.1,  1.03, 7.64),
    @("CB52",   0.12, 1.22, 4.08),
    @("CB101",  0.15, 1.9,  19.38),
    @("CB105",  0.21, 1.36, 10.01),
    @("CB110",  0.14, 1.86, 20.75),
    @("CB118",  0.22, 1.35, 6.1),
    @("CB128",  0.24, 2.34, 26.53),
    @("CB132",  0.12, 1.33, 16.04),
    @("CB138",  0.2,  1.29, 10.85),
    @("CB149",  0.15, 1.35, 14.14),
    @("CB153",  0.22, 1.41, 9.25),
    @("CB156",  0.22, 1.02, 3.85),
    @("CB170",  0.28, 1.16, 4),
    @("CB180",  0.28, 1.62, 7.58)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 3).Value = $data[$i][2]
    $ws.Cells.Item($row, 4).Value = $data[$i][3]
}
